$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 23: CSK vs SRH - fill in the player points for row 32
$ws.Range("E32").Value = 0
$ws.Range("H32").Value = 20
$ws.Range("K32").Value = 40
$ws.Range("N32").Value = 80
$ws.Range("Q32").Value = 100
$ws.Range("T32").Value = 60

$wb.Save()
